# Auto-generated edit script: updates market-price-derived columns (H-N)
# across multiple craft job sheets, per scheduled price refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1928208.8
$ws.Range("I15").Value = 1928208.8
$ws.Range("K15").Value = 5784626.4
$ws.Range("M15").Value = -5784457.4
$ws.Range("H40").Value = 2795.7273
$ws.Range("I40").Value = 2444.5557
$ws.Range("J40").Value = 3038.8462
$ws.Range("K40").Value = 2444.5557
$ws.Range("L40").Value = 3038.8462
$ws.Range("M40").Value = -2269.5557
$ws.Range("N40").Value = -3388.8462
$ws.Range("H64").Value = 3973
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 3970.5454
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 3970.5454
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4466.5454
$ws.Range("H67").Value = 3973
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 3970.5454
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 3970.5454
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5686.5454
$ws.Range("H137").Value = 1259.1333
$ws.Range("I137").Value = 1381.5714
$ws.Range("J137").Value = 973.44446
$ws.Range("K137").Value = 4144.7142
$ws.Range("L137").Value = 2920.33338
$ws.Range("M137").Value = -1594.7142
$ws.Range("N137").Value = -8020.33338

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 37149.855
$ws.Range("I21").Value = 30005
$ws.Range("J21").Value = 42508.5
$ws.Range("K21").Value = 30005
$ws.Range("L21").Value = 42508.5
$ws.Range("M21").Value = -29631
$ws.Range("N21").Value = -43256.5
$ws.Range("H32").Value = 10113.437
$ws.Range("I32").Value = 5225.7
$ws.Range("J32").Value = 30239.412
$ws.Range("K32").Value = 5225.7
$ws.Range("L32").Value = 30239.412
$ws.Range("M32").Value = -4938.7
$ws.Range("N32").Value = -30813.412
$ws.Range("H61").Value = 2128.9167
$ws.Range("I61").Value = 2303.7144
$ws.Range("J61").Value = 1884.2
$ws.Range("K61").Value = 2303.7144
$ws.Range("L61").Value = 1884.2
$ws.Range("M61").Value = -2091.7144
$ws.Range("N61").Value = -2308.2
$ws.Range("H92").Value = 45150
$ws.Range("J92").Value = 45150
$ws.Range("L92").Value = 45150
$ws.Range("N92").Value = -50142
$ws.Range("H109").Value = 16059
$ws.Range("J109").Value = 16059
$ws.Range("L109").Value = 16059
$ws.Range("N109").Value = -18833
$ws.Range("H136").Value = 2128.9167
$ws.Range("I136").Value = 2303.7144
$ws.Range("J136").Value = 1884.2
$ws.Range("K136").Value = 6911.1432
$ws.Range("L136").Value = 5652.6
$ws.Range("M136").Value = -4361.1432
$ws.Range("N136").Value = -10752.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 366.66666
$ws.Range("I44").Value = 366.66666
$ws.Range("K44").Value = 366.66666
$ws.Range("M44").Value = 75.33334000000002
$ws.Range("H132").Value = 3088.3635
$ws.Range("I132").Value = 3361.25
$ws.Range("J132").Value = 3001.04
$ws.Range("K132").Value = 10083.75
$ws.Range("L132").Value = 9003.119999999999
$ws.Range("M132").Value = -7553.75
$ws.Range("N132").Value = -14063.12
$ws.Range("H134").Value = 1494.5763
$ws.Range("I134").Value = 793.8857400000001
$ws.Range("J134").Value = 2516.4167
$ws.Range("K134").Value = 2381.65722
$ws.Range("L134").Value = 7549.250100000001
$ws.Range("M134").Value = 153.3427799999999
$ws.Range("N134").Value = -12619.2501

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 89.36364
$ws.Range("I40").Value = 93.3
$ws.Range("J40").Value = 50
$ws.Range("K40").Value = 373.2
$ws.Range("L40").Value = 200
$ws.Range("M40").Value = -304.2
$ws.Range("N40").Value = -338
$ws.Range("H41").Value = 1930.9565
$ws.Range("I41").Value = 1835.3334
$ws.Range("J41").Value = 1945.3
$ws.Range("K41").Value = 5506.0002
$ws.Range("L41").Value = 5835.9
$ws.Range("M41").Value = -5168.0002
$ws.Range("N41").Value = -6511.9
$ws.Range("H43").Value = 4866.6665
$ws.Range("J43").Value = 4866.6665
$ws.Range("L43").Value = 14599.9995
$ws.Range("N43").Value = -14827.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4278055.5
$ws.Range("I14").Value = 5500071.5
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 5500071.5
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = -5499903.5
$ws.Range("N14").Value = -1336
$ws.Range("H31").Value = 858.3
$ws.Range("I31").Value = 858.3
$ws.Range("K31").Value = 858.3
$ws.Range("M31").Value = -566.3
$ws.Range("H37").Value = 858.3
$ws.Range("I37").Value = 858.3
$ws.Range("K37").Value = 858.3
$ws.Range("M37").Value = -581.3
$ws.Range("H132").Value = 788430.4399999999
$ws.Range("I132").Value = 1345634.6
$ws.Range("K132").Value = 4036903.8
$ws.Range("M132").Value = -4034373.8
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 64166
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 64166
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 64166
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -74526

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9092272
$ws.Range("I136").Value = 12500699
$ws.Range("J136").Value = 3133.3333
$ws.Range("K136").Value = 37502097
$ws.Range("L136").Value = 9399.999899999999
$ws.Range("M136").Value = -37499547
$ws.Range("N136").Value = -14499.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5846
$ws.Range("I62").Value = 4654
$ws.Range("J62").Value = 7038
$ws.Range("K62").Value = 4654
$ws.Range("L62").Value = 7038
$ws.Range("M62").Value = -4030
$ws.Range("N62").Value = -8286
$ws.Range("H65").Value = 5846
$ws.Range("I65").Value = 4654
$ws.Range("J65").Value = 7038
$ws.Range("K65").Value = 23270
$ws.Range("L65").Value = 35190
$ws.Range("M65").Value = -20150
$ws.Range("N65").Value = -41430
$ws.Range("H132").Value = 1608.0273
$ws.Range("I132").Value = 1438.5853
$ws.Range("J132").Value = 1825.125
$ws.Range("K132").Value = 4315.7559
$ws.Range("L132").Value = 5475.375
$ws.Range("M132").Value = -1785.7559
$ws.Range("N132").Value = -10535.375
$ws.Range("H136").Value = 4722714
$ws.Range("I136").Value = 5958848.5
$ws.Range("K136").Value = 17876545.5
$ws.Range("M136").Value = -17873995.5
